$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 787
$ws.Range("F4").Value = 36
$ws.Range("F5").Value = 225
$ws.Range("F7").Value = 1144
$ws.Range("F8").Value = 899
$ws.Range("F10").Value = 716
$ws.Range("F12").Value = 1425
$ws.Range("F15").Value = 1594
$ws.Range("F16").Value = 14
$ws.Range("F17").Value = 603
$ws.Range("F18").Value = 16
$ws.Range("F20").Value = 371
$ws.Range("F23").Value = 748
$ws.Range("F30").Value = 297
$ws.Range("F31").Value = 2413
$ws.Range("F32").Value = 277
$ws.Range("F33").Value = 1350
$ws.Range("F34").Value = 458
$ws.Range("F36").Value = 3936

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F14").Value = 4135
$ws.Range("F20").Value = 254
$ws.Range("F28").Value = 1710

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 1267
$ws.Range("F5").Value = 1661
$ws.Range("F6").Value = 451
$ws.Range("F7").Value = 994

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1267
$ws.Range("F4").Value = 1661
$ws.Range("F5").Value = 451
$ws.Range("F6").Value = 994
$ws.Range("F8").Value = 787
$ws.Range("F9").Value = 36
$ws.Range("F10").Value = 225
$ws.Range("F12").Value = 1144
$ws.Range("F13").Value = 899
$ws.Range("F17").Value = 716
$ws.Range("F22").Value = 1425
$ws.Range("F25").Value = 1594
$ws.Range("F26").Value = 14
$ws.Range("F27").Value = 603
$ws.Range("F29").Value = 371
$ws.Range("F32").Value = 748
$ws.Range("F38").Value = 254
$ws.Range("F42").Value = 297
$ws.Range("F43").Value = 2413
$ws.Range("F46").Value = 1710
$ws.Range("F47").Value = 1710
$ws.Range("F48").Value = 1350
$ws.Range("F49").Value = 458
$ws.Range("F50").Value = 3936
